# Add equations to t0, tg and tpeak definitions on Sheet1.
#
# Row 3  (B3 = FirstGermTime) -> C3: add the t0 equation block
# Row 4  (B4 = LastGermTime)  -> C4: add the tg equation block
# Row 6  (B6 = PeakGermTime)  -> C6: add the tpeak equation block
#                              -> A6: append "($t_{peak}$)" to the label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$c3Text = @'
It is the time for first germination to occur (e.g. First day of germination).
$$t_{0} = \min \left\{ T_{i} : N_{i} \neq  0 \right\}$$
Where, $T_{i}$ is the time from the start of the experiment to the $i$th interval and $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval)
'@

$c4Text = @'
It is the time for last germination to occur (e.g. Last day of germination)
$$t_{g} = \max \left\{ T_{i} : N_{i} \neq  0 \right\}$$
Where, $T_{i}$ is the time from the start of the experiment to the $i$th interval and $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval)
'@

$c6Text = @'
It is the time in which highest frequency of germinated seeds are observed and need not be unique.
$$t_{peak} = \left\{ T_{i} : N_{i} = N_{max} \right\}$$
Where, $T_{i}$ is the time from the start of the experiment to the $i$th interval, $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval) and $N_{max}$ is the maximum number of seeds germinated per interval.
'@

$a6Text = 'Peak period of germination or Modal time of germination ($t_{peak}$)'

# Shared-string table append order matters for an exact byte match: the
# author's commit ends up with the "last germination" (tg) text before the
# "first germination" (t0) text, so write C4 ahead of C3.
$ws.Range("C4").Value = $c4Text
$ws.Range("C3").Value = $c3Text
$ws.Range("C6").Value = $c6Text
$ws.Range("A6").Value = $a6Text

# Mirror the author's final view state: A6 selected (last-edited cell).
$ws.Range("A6").Select()
